$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two rows "Altair AcuSolve" (row 5) and "PACEFISH" (row 6) are removed;
# everything below shifts up by two rows.
$ws.Range("A5:A6").EntireRow.Delete()

# Fill in new chapter-2 content for the rows that used to hold
# HIFUN / DYVERSO / ANSYS FLUENT / NUMERIX (now rows 5-8).
$ws.Range("C5").Value = "复杂几何外形的高升力流动问题"
$ws.Range("D5").Value = "单节点、多GPU"
$ws.Range("D5").WrapText = $false

$ws.Range("B6").Value = "SPH"
$ws.Rows.Item(6).RowHeight = 25

$ws.Range("B7").Value = "结构化/非结构化"
$ws.Range("C7").Value = "通用CFD商业软件"
$ws.Rows.Item(7).RowHeight = 25

$ws.Range("B8").Value = "SPH"
$ws.Rows.Item(8).RowHeight = 25

$ws.Range("D8").Select()
